$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:C2,F2,A3:F3,G3,A4:B4,E4:F4").HorizontalAlignment = -4142
